# Applies the scheduled-runner profit recalculation update to all 8 leve sheets.
# For each affected row, the cost & profit columns (H,I,J,K,L,M,N) are refreshed
# to the newly computed values; a few rows gain or lose their trailing M/N cell
# entirely where the recompute no longer produces (or now produces) a value.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3666.5557
$ws.Range("J17").Value = 3666.5557
$ws.Range("L17").Value = 10999.6671
$ws.Range("N17").Value = -11335.6671
$ws.Range("H74").Value = 8000
$ws.Range("I74").Value = 8000
$ws.Range("K74").Value = 8000
$ws.Range("M74").Value = -7064
$ws.Range("H77").Value = 8000
$ws.Range("I77").Value = 8000
$ws.Range("K77").Value = 40000
$ws.Range("M77").Value = -35320
$ws.Range("H92").Value = 420.625
$ws.Range("I92").Value = 422.5
$ws.Range("K92").Value = 422.5
$ws.Range("M92").Value = 825.5
$ws.Range("H98").Value = 1225.5834
$ws.Range("I98").Value = 1225.5834
$ws.Range("K98").Value = 1225.5834
$ws.Range("M98").Value = 272.4166
$ws.Range("H100").Value = 1127
$ws.Range("I100").Value = 776.05
$ws.Range("K100").Value = 776.05
$ws.Range("M100").Value = -235.05
$ws.Range("H121").Value = 6736.75
$ws.Range("J121").Value = 6736.75
$ws.Range("L121").Value = 20210.25
$ws.Range("N121").Value = -23704.25
$ws.Range("H122").Value = 1225.5834
$ws.Range("I122").Value = 1225.5834
$ws.Range("K122").Value = 3676.7502
$ws.Range("M122").Value = -1226.7502
$ws.Range("H125").Value = 1238.6
$ws.Range("I125").Value = 1238.6
$ws.Range("K125").Value = 11147.4
$ws.Range("M125").Value = -8687.4
$ws.Range("H137").Value = 26425
$ws.Range("J137").Value = 29995.834
$ws.Range("L137").Value = 89987.50199999999
$ws.Range("N137").Value = -95087.50199999999
$ws.Range("H138").Value = 7278.1606
$ws.Range("I138").Value = 7365.636
$ws.Range("K138").Value = 22096.908
$ws.Range("M138").Value = -16956.908

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23384
$ws.Range("I32").Value = 13656.682
$ws.Range("J32").Value = 29868.879
$ws.Range("K32").Value = 13656.682
$ws.Range("L32").Value = 29868.879
$ws.Range("M32").Value = -13369.682
$ws.Range("N32").Value = -30442.879
$ws.Range("H43").Value = 53000
$ws.Range("J43").Value = 53000
$ws.Range("L43").Value = 53000
$ws.Range("N43").Value = -53626
$ws.Range("H95").Value = 54749.5
$ws.Range("J95").Value = 54749.5
$ws.Range("L95").Value = 54749.5
$ws.Range("N95").Value = -60241.5
$ws.Range("H102").Value = 1861.409
$ws.Range("I102").Value = 1469.3
$ws.Range("K102").Value = 1469.3
$ws.Range("M102").Value = 152.7
$ws.Range("H122").Value = 3750.625
$ws.Range("I122").Value = 3635.818
$ws.Range("J122").Value = 3847.7693
$ws.Range("K122").Value = 10907.454
$ws.Range("L122").Value = 11543.3079
$ws.Range("M122").Value = -8457.454000000002
$ws.Range("N122").Value = -16443.3079
$ws.Range("H132").Value = 3276.2354
$ws.Range("I132").Value = 1833.0714
$ws.Range("K132").Value = 5499.2142
$ws.Range("M132").Value = -2969.2142

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H94").Value = 1966.5
$ws.Range("J94").Value = 2503
$ws.Range("L94").Value = 2503
$ws.Range("N94").Value = -3405
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3141.375
$ws.Range("I105").Value = 1603.4445
$ws.Range("K105").Value = 1603.4445
$ws.Range("M105").Value = 143.5554999999999
$ws.Range("H109").Value = 54097.727
$ws.Range("J109").Value = 54097.727
$ws.Range("L109").Value = 54097.727
$ws.Range("N109").Value = -56177.727
$ws.Range("H132").Value = 2899.9524
$ws.Range("I132").Value = 2786.111
$ws.Range("K132").Value = 8358.332999999999
$ws.Range("M132").Value = -5828.332999999999
$ws.Range("H134").Value = 3534.9473
$ws.Range("I134").Value = 2780.7693
$ws.Range("K134").Value = 8342.3079
$ws.Range("M134").Value = -5807.3079
$ws.Range("H141").Value = 121000
$ws.Range("J141").Value = 121000
$ws.Range("L141").Value = 121000
$ws.Range("N141").Value = -131360

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 8296.799999999999
$ws.Range("I59").Value = 7542
$ws.Range("J59").Value = 8800
$ws.Range("K59").Value = 22626
$ws.Range("L59").Value = 26400
$ws.Range("M59").Value = -22086
$ws.Range("N59").Value = -27480
$ws.Range("H113").Value = 5990
$ws.Range("J113").Value = 5990
$ws.Range("L113").Value = 17970
$ws.Range("N113").Value = -22310
$ws.Range("H121").Value = 1207.0834
$ws.Range("I121").Value = 804.5
$ws.Range("J121").Value = 1609.6666
$ws.Range("K121").Value = 2413.5
$ws.Range("L121").Value = 4828.9998
$ws.Range("M121").Value = -1103.5
$ws.Range("N121").Value = -7448.9998
$ws.Range("H122").Value = 2113.7144
$ws.Range("I122").Value = 1449.5
$ws.Range("J122").Value = 2999.3333
$ws.Range("K122").Value = 13045.5
$ws.Range("L122").Value = 26993.9997
$ws.Range("M122").Value = -10595.5
$ws.Range("N122").Value = -31893.9997

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8500.125
$ws.Range("I80").Value = 9995
$ws.Range("J80").Value = 8286.571
$ws.Range("K80").Value = 9995
$ws.Range("L80").Value = 8286.571
$ws.Range("M80").Value = -8997
$ws.Range("N80").Value = -10282.571
$ws.Range("H83").Value = 8500.125
$ws.Range("I83").Value = 9995
$ws.Range("J83").Value = 8286.571
$ws.Range("K83").Value = 49975
$ws.Range("L83").Value = 41432.855
$ws.Range("M83").Value = -44983
$ws.Range("N83").Value = -51416.855
$ws.Range("H102").Value = 4035.3333
$ws.Range("I102").Value = 2829.5
$ws.Range("K102").Value = 2829.5
$ws.Range("M102").Value = -1207.5
$ws.Range("H122").Value = 924064.9399999999
$ws.Range("I122").Value = 210158.2
$ws.Range("K122").Value = 630474.6000000001
$ws.Range("M122").Value = -628024.6000000001
$ws.Range("H132").Value = 3487.3044
$ws.Range("I132").Value = 2677.25
$ws.Range("J132").Value = 8887.666999999999
$ws.Range("K132").Value = 8031.75
$ws.Range("L132").Value = 26663.001
$ws.Range("M132").Value = -5501.75
$ws.Range("N132").Value = -31723.001

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2600.5
$ws.Range("H27").Value = 2600.5
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H46").Value = 3250
$ws.Range("I46").Value = 2714.2856
$ws.Range("J46").Value = 3666.6667
$ws.Range("K46").Value = 2714.2856
$ws.Range("L46").Value = 3666.6667
$ws.Range("M46").Value = -2526.2856
$ws.Range("N46").Value = -4042.6667
$ws.Range("H55").Value = 1146.8948
$ws.Range("I55").Value = 986.8125
$ws.Range("J55").Value = 2000.6666
$ws.Range("K55").Value = 986.8125
$ws.Range("L55").Value = 2000.6666
$ws.Range("M55").Value = -813.8125
$ws.Range("N55").Value = -2346.6666
$ws.Range("H61").Value = 5550.24
$ws.Range("I61").Value = 5633.2383
$ws.Range("K61").Value = 5633.2383
$ws.Range("M61").Value = -5431.2383
$ws.Range("H100").Value = 3361.6365
$ws.Range("I100").Value = 3222.5557
$ws.Range("J100").Value = 3987.5
$ws.Range("K100").Value = 3222.5557
$ws.Range("L100").Value = 3987.5
$ws.Range("M100").Value = -2681.5557
$ws.Range("N100").Value = -5069.5
$ws.Range("H110").Value = 59548
$ws.Range("J110").Value = 59548
$ws.Range("L110").Value = 59548
$ws.Range("N110").Value = -67728
$ws.Range("H113").Value = 5550.24
$ws.Range("I113").Value = 5633.2383
$ws.Range("K113").Value = 5633.2383
$ws.Range("M113").Value = -3463.2383
$ws.Range("H132").Value = 5911.4287
$ws.Range("I132").Value = 2380
$ws.Range("K132").Value = 7140
$ws.Range("M132").Value = -4610
$ws.Range("H136").Value = 3994.4
$ws.Range("I136").Value = 3993.25
$ws.Range("K136").Value = 11979.75
$ws.Range("M136").Value = -9429.75

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 11967.143
$ws.Range("J5").Value = 11967.143
$ws.Range("L5").Value = 11967.143
$ws.Range("N5").Value = -12191.143
$ws.Range("H28").Value = 30017.2
$ws.Range("I28").Value = 30014
$ws.Range("K28").Value = 30014
$ws.Range("M28").Value = -29666
$ws.Range("H32").Value = 247900
$ws.Range("I32").Value = 363166.66
$ws.Range("J32").Value = 75000
$ws.Range("K32").Value = 363166.66
$ws.Range("L32").Value = 75000
$ws.Range("M32").Value = -362849.66
$ws.Range("N32").Value = -75634
$ws.Range("H113").Value = 1240.0588
$ws.Range("I113").Value = 583
$ws.Range("K113").Value = 1749
$ws.Range("M113").Value = 421
$ws.Range("H122").Value = 1747.7858
$ws.Range("I122").Value = 1747.7858
$ws.Range("K122").Value = 5243.357400000001
$ws.Range("M122").Value = -2793.357400000001
$ws.Range("H132").Value = 2312.5334
$ws.Range("I132").Value = 1188.1111
$ws.Range("K132").Value = 3564.3333
$ws.Range("M132").Value = -1034.3333
